# edit.ps1
# Applies the CS:GO "CT-side" document edit:
#  - Re-splits runs around English/gaming loanwords and wraps them in
#    <w:proofErr w:type="spellStart"/> ... <w:proofErr w:type="spellEnd"/>
#    markers, exactly as Word's proofing pass does when it re-saves a
#    document that mixes Dutch body text with English terminology.
#  - Normalizes paragraph properties (the explicit pPr/pStyle blocks
#    collapse away because they only ever dictated the built-in
#    "Normal" style with no overrides).
#  - Adds a <w:lastRenderedPageBreak/> before item 6, matching the
#    pagination Word recorded on save.

$d = $word.ActiveDocument

$paras = @(
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:rPr>
          <w:sz w:val="48"/>
          <w:szCs w:val="48"/>
        </w:rPr>
        <w:t>CT-side</w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">In Counter-Strike: Global </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Offensive</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> (CS:GO) verwijst &quot;CT-side&quot; naar de zijde van het spel waarin je speelt als de Counter-</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Terrorists</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> (</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>CTs</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">). De </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>CTs</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> zijn verantwoordelijk voor het beschermen van specifieke doelen, zoals het voorkomen van bomaanslagen (</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>defusals</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>) of gijzelaars (</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>hostages</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>) redden, afhankelijk van de spelmodus. Hier is een uitleg over de CT-side in CS:GO:</w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t>1. **Doel van de CT-side**:</w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">   - Het belangrijkste doel van de CT-side is om de aanvallende </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Terrorists</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> (T''s) te stoppen. Dit betekent meestal het voorkomen van bomplantages in de bomplaatsen of het beveiligen van gijzelaars om een gijzeling te winnen. Het uiteindelijke doel is om rondes te winnen door ofwel alle tegenstanders te elimineren, de bom te ontmantelen, of gijzelaars te redden.</w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t>2. **Uitrusting**:</w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">   - Als CT heb je toegang tot verschillende wapens en uitrusting, zoals </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>rifles</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>submachine</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>guns</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>shotguns</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>, pistolen, granaten en pantser. Het is essentieel om je geld verstandig te beheren en de juiste uitrusting te kopen, afhankelijk van je rol en de economische situatie van je team.</w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t>3. **Verdedigen van doelen**:</w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">   - Als CT is het belangrijk om bomplaatsen en </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>gijzelaarlocaties</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> te beveiligen. Dit vereist communicatie en samenwerking met je team om strategische posities in te nemen en de aanval van de T''s af te weren. Het bewaken van chokepoints en het gebruik van granaten om vijanden te vertragen of te verzwakken, is van cruciaal belang.</w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t>4. **Economisch beheer**:</w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">   - Als CT is het essentieel om je economie goed te beheren. Wanneer je verliest, verdien je minder geld en moet je beslissen of je wapens en uitrusting wilt behouden voor de volgende ronde of wilt kopen met beperkte middelen.</w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t>5. **Rotaties en communicatie**:</w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">   - Communicatie is van het grootste belang op de CT-side. Spelers moeten informatie delen over de positie van de vijand en beslissingen nemen over wanneer en hoe ze moeten roteren om de doelen te beschermen. Te late of ongecoördineerde rotaties kunnen het verlies van een ronde tot gevolg hebben.</w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>6. **Het spelen van de klok**:</w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">   - Als </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>CT''s</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> kun je soms proberen de tijd op de </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>rondeklok</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> te spelen. Door de T''s onder druk te zetten en de klok te laten wegtikken, kun je hen dwingen overhaaste beslissingen te nemen die hun slagingskansen verminderen.</w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>',
'<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:r>
        <w:t xml:space="preserve">De CT-side vereist een goede teamcoördinatie, communicatie en een begrip van de </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>maplay</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>-out om effectief te zijn. Het is van vitaal belang om samen te werken met je teamgenoten en om strategisch te denken om te voorkomen dat de T''s hun doel bereiken en om ronden te winnen.</w:t>
      </w:r>
    </w:p>'
)

for ($i = 0; $i -lt $paras.Count; $i++) {
    $d.Paragraphs($i + 1).Range.InsertXML($paras[$i])
}

# The page stays portrait-oriented; re-assert it so PageSetup stays in
# sync with the (unchanged) pgSz dimensions.
$d.PageSetup.Orientation = 0

Write-Host ("paragraphs: " + $d.Paragraphs.Count)
